$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 79.14286
$ws.Range("I9").Value = 79.14286
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 79.14286
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = 89.85714
# Row 40
$ws.Range("H40").Value = 3000
$ws.Range("J40").Value = 3000
$ws.Range("L40").Value = 3000
$ws.Range("N40").Value = -3350
# Row 86
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = 0
# Row 89
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = 0
# Row 120
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").ClearContents()
$ws.Range("N120").Value = 0
# Row 138
$ws.Range("H138").Value = 2049.53
$ws.Range("I138").Value = 1403.4783
$ws.Range("J138").Value = 2242.5066
$ws.Range("K138").Value = 4210.4349
$ws.Range("L138").Value = 6727.5198
$ws.Range("M138").Value = 929.5650999999998
$ws.Range("N138").Value = -17007.5198

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 327725.06
$ws.Range("J2").Value = 1089.8889
$ws.Range("L2").Value = 1089.8889
$ws.Range("N2").Value = -1315.8889
# Row 45
$ws.Range("H45").Value = 9001631
$ws.Range("I45").Value = 22501228
$ws.Range("K45").Value = 22501228
$ws.Range("M45").Value = -22500851
# Row 63
$ws.Range("H63").Value = 2054
$ws.Range("J63").Value = 2651.5715
$ws.Range("L63").Value = 2651.5715
$ws.Range("N63").Value = -4023.5715
# Row 66
$ws.Range("H66").Value = 2054
$ws.Range("J66").Value = 2651.5715
$ws.Range("L66").Value = 13257.8575
$ws.Range("N66").Value = -20121.8575
# Row 110
$ws.Range("H110").Value = 1191.7273
$ws.Range("I110").Value = 862
$ws.Range("J110").Value = 2071
$ws.Range("K110").Value = 862
$ws.Range("L110").Value = 2071
$ws.Range("M110").Value = 1183
$ws.Range("N110").Value = -6161
# Row 116
$ws.Range("H116").Value = 327725.06
$ws.Range("J116").Value = 1089.8889
$ws.Range("L116").Value = 1089.8889
$ws.Range("N116").Value = -5677.8889
# Row 122
$ws.Range("H122").Value = 1277
$ws.Range("I122").Value = 1248.0333
$ws.Range("J122").Value = 1566.6666
$ws.Range("K122").Value = 3744.0999
$ws.Range("L122").Value = 4699.9998
$ws.Range("M122").Value = -1294.0999
$ws.Range("N122").Value = -9599.9998

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 327725.06
$ws.Range("J3").Value = 1089.8889
$ws.Range("L3").Value = 1089.8889
$ws.Range("N3").Value = -1317.8889
# Row 82
$ws.Range("H82").Value = 21101
$ws.Range("J82").Value = 32400
$ws.Range("L82").Value = 32400
$ws.Range("N82").Value = -33166
# Row 85
$ws.Range("H85").Value = 21101
$ws.Range("J85").Value = 32400
$ws.Range("L85").Value = 32400
$ws.Range("N85").Value = -35052
# Row 86
$ws.Range("H86").Value = 1000874.5
$ws.Range("I86").Value = 1749
$ws.Range("K86").Value = 1749
$ws.Range("M86").Value = -626
# Row 89
$ws.Range("H89").Value = 1000874.5
$ws.Range("I89").Value = 1749
$ws.Range("K89").Value = 8745
$ws.Range("M89").Value = -3129
# Row 105
$ws.Range("H105").Value = 2155.3242
$ws.Range("I105").Value = 2148.4688
$ws.Range("J105").Value = 2199.2
$ws.Range("K105").Value = 2148.4688
$ws.Range("L105").Value = 2199.2
$ws.Range("M105").Value = -401.4688000000001
$ws.Range("N105").Value = -5693.2
# Row 107
$ws.Range("H107").Value = 2414.3333
$ws.Range("J107").Value = 1856.5
$ws.Range("L107").Value = 1856.5
$ws.Range("N107").Value = -5696.5
# Row 134
$ws.Range("H134").Value = 3482.1667
$ws.Range("I134").Value = 3175.0977
$ws.Range("J134").Value = 5280.7144
$ws.Range("K134").Value = 9525.293099999999
$ws.Range("L134").Value = 15842.1432
$ws.Range("M134").Value = -6990.293099999999
$ws.Range("N134").Value = -20912.1432

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1188.25
$ws.Range("I16").Value = 1188.25
$ws.Range("K16").Value = 1188.25
$ws.Range("M16").Value = -901.25
# Row 31
$ws.Range("H31").Value = 2289.963
$ws.Range("I31").Value = 1646.1818
$ws.Range("J31").Value = 2732.5625
$ws.Range("K31").Value = 1646.1818
$ws.Range("L31").Value = 2732.5625
$ws.Range("M31").Value = -1351.1818
$ws.Range("N31").Value = -3322.5625
# Row 34
$ws.Range("H34").Value = 2289.963
$ws.Range("I34").Value = 1646.1818
$ws.Range("J34").Value = 2732.5625
$ws.Range("K34").Value = 1646.1818
$ws.Range("L34").Value = 2732.5625
$ws.Range("M34").Value = -1444.1818
$ws.Range("N34").Value = -3136.5625
# Row 113
$ws.Range("H113").Value = 1188.25
$ws.Range("I113").Value = 1188.25
$ws.Range("K113").Value = 1188.25
$ws.Range("M113").Value = 981.75
# Row 132
$ws.Range("H132").Value = 1371.3684
$ws.Range("I132").Value = 934.40625
$ws.Range("J132").Value = 3701.8333
$ws.Range("K132").Value = 2803.21875
$ws.Range("L132").Value = 11105.4999
$ws.Range("M132").Value = -273.21875
$ws.Range("N132").Value = -16165.4999

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 126
$ws.Range("H126").Value = 4750
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 5833.3335
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 17500.0005
$ws.Range("M126").Value = 440
$ws.Range("N126").Value = -27380.0005
# Row 131
$ws.Range("H131").Value = 14258.328
$ws.Range("I131").Value = 490.75
$ws.Range("J131").Value = 15224.474
$ws.Range("K131").Value = 1472.25
$ws.Range("L131").Value = 45673.422
$ws.Range("M131").Value = 3567.75
$ws.Range("N131").Value = -55753.422

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4842
$ws.Range("I70").Value = 4700
$ws.Range("J70").Value = 4877.5
$ws.Range("K70").Value = 4700
$ws.Range("L70").Value = 4877.5
$ws.Range("M70").Value = -4430
$ws.Range("N70").Value = -5417.5
# Row 73
$ws.Range("H73").Value = 4842
$ws.Range("I73").Value = 4700
$ws.Range("J73").Value = 4877.5
$ws.Range("K73").Value = 4700
$ws.Range("L73").Value = 4877.5
$ws.Range("M73").Value = -3764
$ws.Range("N73").Value = -6749.5
# Row 97
$ws.Range("H97").Value = 1182.5385
$ws.Range("I97").Value = 1060.579
$ws.Range("K97").Value = 1060.579
$ws.Range("M97").Value = -564.579
# Row 113
$ws.Range("H113").Value = 1366.5
$ws.Range("I113").Value = 700
$ws.Range("K113").Value = 700
$ws.Range("M113").Value = 1470
# Row 122
$ws.Range("H122").Value = 958.5833
$ws.Range("I122").Value = 974.0909
$ws.Range("K122").Value = 2922.2727
$ws.Range("M122").Value = -472.2727

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 9478.857
$ws.Range("I16").Value = 9478.857
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 9478.857
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -9308.857
# Row 93
$ws.Range("H93").Value = 1749.3572
$ws.Range("I93").Value = 984.0909
$ws.Range("J93").Value = 4555.3335
$ws.Range("K93").Value = 984.0909
$ws.Range("L93").Value = 4555.3335
$ws.Range("M93").Value = 263.9091
$ws.Range("N93").Value = -7051.3335
# Row 100
$ws.Range("H100").Value = 2002.1818
$ws.Range("I100").Value = 1103.5714
$ws.Range("K100").Value = 1103.5714
$ws.Range("M100").Value = -562.5714

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 600.63635
$ws.Range("I107").Value = 460.75
$ws.Range("J107").Value = 1999.5
$ws.Range("K107").Value = 1382.25
$ws.Range("L107").Value = 5998.5
$ws.Range("M107").Value = 537.75
$ws.Range("N107").Value = -9838.5
# Row 126
$ws.Range("H126").Value = 7960.278
$ws.Range("I126").Value = 10369.167
$ws.Range("J126").Value = 3142.5
$ws.Range("K126").Value = 31107.501
$ws.Range("L126").Value = 9427.5
$ws.Range("M126").Value = -28637.501
$ws.Range("N126").Value = -14367.5
